$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 93.32574466666667
$ws.Cells.Item(2, 8).Value = 279.977234
$ws.Cells.Item(2, 9).Value = 0.2327963689879921
$ws.Cells.Item(2, 10).Value = 0.2327963689879922
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 13.441269
$ws.Cells.Item(2, 14).Value = 40.323807
$ws.Cells.Item(2, 15).Value = 0.08973082133481231
$ws.Cells.Item(2, 16).Value = 0.08973082133481232
$ws.Cells.Item(2, 17).Value = 1254.416438689982
$ws.Cells.Item(2, 18).Value = 11289.74794820984
$ws.Cells.Item(2, 19).Value = 0.02088900939305456
$ws.Cells.Item(2, 20).Value = 0.02088900939305457
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 93.32574466666667
$ws.Cells.Item(3, 8).Value = 279.977234
$ws.Cells.Item(3, 9).Value = 0.2327963689879921
$ws.Cells.Item(3, 10).Value = 0.2327963689879922
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 54.711535
$ws.Cells.Item(3, 14).Value = 164.134605
$ws.Cells.Item(3, 15).Value = 0.3652416280068742
$ws.Cells.Item(3, 16).Value = 0.3652416280068742
$ws.Cells.Item(3, 17).Value = 5105.994745731397
$ws.Cells.Item(3, 18).Value = 45953.95271158257
$ws.Cells.Item(3, 19).Value = 0.08502692480326322
$ws.Cells.Item(3, 20).Value = 0.08502692480326327
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 93.32574466666667
$ws.Cells.Item(4, 8).Value = 279.977234
$ws.Cells.Item(4, 9).Value = 0.2327963689879921
$ws.Cells.Item(4, 10).Value = 0.2327963689879922
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 63.67711
$ws.Cells.Item(4, 14).Value = 191.03133
$ws.Cells.Item(4, 15).Value = 0.4250937452800914
$ws.Cells.Item(4, 16).Value = 0.4250937452800915
$ws.Cells.Item(4, 17).Value = 5942.713708971247
$ws.Cells.Item(4, 18).Value = 53484.42338074122
$ws.Cells.Item(4, 19).Value = 0.0989602803807117
$ws.Cells.Item(4, 20).Value = 0.09896028038071172
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 93.32574466666667
$ws.Cells.Item(5, 8).Value = 279.977234
$ws.Cells.Item(5, 9).Value = 0.2327963689879921
$ws.Cells.Item(5, 10).Value = 0.2327963689879922
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 17.96553866666667
$ws.Cells.Item(5, 14).Value = 53.896616
$ws.Cells.Item(5, 15).Value = 0.119933805378222
$ws.Cells.Item(5, 16).Value = 0.119933805378222
$ws.Cells.Item(5, 17).Value = 1676.64727440446
$ws.Cells.Item(5, 18).Value = 15089.82546964015
$ws.Cells.Item(5, 19).Value = 0.0279201544109626
$ws.Cells.Item(5, 20).Value = 0.02792015441096261
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 164.7897643333334
$ws.Cells.Item(6, 8).Value = 494.369293
$ws.Cells.Item(6, 9).Value = 0.4110597662007076
$ws.Cells.Item(6, 10).Value = 0.4110597662007077
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 13.441269
$ws.Cells.Item(6, 14).Value = 40.323807
$ws.Cells.Item(6, 15).Value = 0.08973082133481231
$ws.Cells.Item(6, 16).Value = 0.08973082133481232
$ws.Cells.Item(6, 17).Value = 2214.983550850939
$ws.Cells.Item(6, 18).Value = 19934.85195765845
$ws.Cells.Item(6, 19).Value = 0.03688473043888541
$ws.Cells.Item(6, 20).Value = 0.03688473043888543
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 164.7897643333334
$ws.Cells.Item(7, 8).Value = 494.369293
$ws.Cells.Item(7, 9).Value = 0.4110597662007076
$ws.Cells.Item(7, 10).Value = 0.4110597662007077
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 54.711535
$ws.Cells.Item(7, 14).Value = 164.134605
$ws.Cells.Item(7, 15).Value = 0.3652416280068742
$ws.Cells.Item(7, 16).Value = 0.3652416280068742
$ws.Cells.Item(7, 17).Value = 9015.900958964919
$ws.Cells.Item(7, 18).Value = 81143.10863068426
$ws.Cells.Item(7, 19).Value = 0.1501361382152715
$ws.Cells.Item(7, 20).Value = 0.1501361382152716
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 164.7897643333334
$ws.Cells.Item(8, 8).Value = 494.369293
$ws.Cells.Item(8, 9).Value = 0.4110597662007076
$ws.Cells.Item(8, 10).Value = 0.4110597662007077
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 63.67711
$ws.Cells.Item(8, 14).Value = 191.03133
$ws.Cells.Item(8, 15).Value = 0.4250937452800914
$ws.Cells.Item(8, 16).Value = 0.4250937452800915
$ws.Cells.Item(8, 17).Value = 10493.33595032774
$ws.Cells.Item(8, 18).Value = 94440.0235529497
$ws.Cells.Item(8, 19).Value = 0.1747389355482175
$ws.Cells.Item(8, 20).Value = 0.1747389355482176
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 164.7897643333334
$ws.Cells.Item(9, 8).Value = 494.369293
$ws.Cells.Item(9, 9).Value = 0.4110597662007076
$ws.Cells.Item(9, 10).Value = 0.4110597662007077
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 17.96553866666667
$ws.Cells.Item(9, 14).Value = 53.896616
$ws.Cells.Item(9, 15).Value = 0.119933805378222
$ws.Cells.Item(9, 16).Value = 0.119933805378222
$ws.Cells.Item(9, 17).Value = 2960.536883001388
$ws.Cells.Item(9, 18).Value = 26644.83194701249
$ws.Cells.Item(9, 19).Value = 0.04929996199833311
$ws.Cells.Item(9, 20).Value = 0.04929996199833311
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 95.835818
$ws.Cells.Item(10, 8).Value = 287.507454
$ws.Cells.Item(10, 9).Value = 0.2390576204784642
$ws.Cells.Item(10, 10).Value = 0.2390576204784643
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 13.441269
$ws.Cells.Item(10, 14).Value = 40.323807
$ws.Cells.Item(10, 15).Value = 0.08973082133481231
$ws.Cells.Item(10, 16).Value = 0.08973082133481232
$ws.Cells.Item(10, 17).Value = 1288.155009573042
$ws.Cells.Item(10, 18).Value = 11593.39508615738
$ws.Cells.Item(10, 19).Value = 0.02145083663187844
$ws.Cells.Item(10, 20).Value = 0.02145083663187845
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 95.835818
$ws.Cells.Item(11, 8).Value = 287.507454
$ws.Cells.Item(11, 9).Value = 0.2390576204784642
$ws.Cells.Item(11, 10).Value = 0.2390576204784643
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 54.711535
$ws.Cells.Item(11, 14).Value = 164.134605
$ws.Cells.Item(11, 15).Value = 0.3652416280068742
$ws.Cells.Item(11, 16).Value = 0.3652416280068742
$ws.Cells.Item(11, 17).Value = 5243.32471076063
$ws.Cells.Item(11, 18).Value = 47189.92239684566
$ws.Cells.Item(11, 19).Value = 0.08731379449100374
$ws.Cells.Item(11, 20).Value = 0.08731379449100377
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 95.835818
$ws.Cells.Item(12, 8).Value = 287.507454
$ws.Cells.Item(12, 9).Value = 0.2390576204784642
$ws.Cells.Item(12, 10).Value = 0.2390576204784643
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 63.67711
$ws.Cells.Item(12, 14).Value = 191.03133
$ws.Cells.Item(12, 15).Value = 0.4250937452800914
$ws.Cells.Item(12, 16).Value = 0.4250937452800915
$ws.Cells.Item(12, 17).Value = 6102.54792472598
$ws.Cells.Item(12, 18).Value = 54922.93132253382
$ws.Cells.Item(12, 19).Value = 0.101621899226937
$ws.Cells.Item(12, 20).Value = 0.1016218992269371
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 95.835818
$ws.Cells.Item(13, 8).Value = 287.507454
$ws.Cells.Item(13, 9).Value = 0.2390576204784642
$ws.Cells.Item(13, 10).Value = 0.2390576204784643
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 17.96553866666667
$ws.Cells.Item(13, 14).Value = 53.896616
$ws.Cells.Item(13, 15).Value = 0.119933805378222
$ws.Cells.Item(13, 16).Value = 0.119933805378222
$ws.Cells.Item(13, 17).Value = 1721.742093930629
$ws.Cells.Item(13, 18).Value = 15495.67884537566
$ws.Cells.Item(13, 19).Value = 0.02867109012864499
$ws.Cells.Item(13, 20).Value = 0.02867109012864499
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 46.93870866666666
$ws.Cells.Item(14, 8).Value = 140.816126
$ws.Cells.Item(14, 9).Value = 0.117086244332836
$ws.Cells.Item(14, 10).Value = 0.117086244332836
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 13.441269
$ws.Cells.Item(14, 14).Value = 40.323807
$ws.Cells.Item(14, 15).Value = 0.08973082133481231
$ws.Cells.Item(14, 16).Value = 0.08973082133481232
$ws.Cells.Item(14, 17).Value = 630.9158097012979
$ws.Cells.Item(14, 18).Value = 5678.242287311682
$ws.Cells.Item(14, 19).Value = 0.01050624487099388
$ws.Cells.Item(14, 20).Value = 0.01050624487099389
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 46.93870866666666
$ws.Cells.Item(15, 8).Value = 140.816126
$ws.Cells.Item(15, 9).Value = 0.117086244332836
$ws.Cells.Item(15, 10).Value = 0.117086244332836
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 54.711535
$ws.Cells.Item(15, 14).Value = 164.134605
$ws.Cells.Item(15, 15).Value = 0.3652416280068742
$ws.Cells.Item(15, 16).Value = 0.3652416280068742
$ws.Cells.Item(15, 17).Value = 2568.088802071136
$ws.Cells.Item(15, 18).Value = 23112.79921864023
$ws.Cells.Item(15, 19).Value = 0.04276477049733565
$ws.Cells.Item(15, 20).Value = 0.04276477049733567
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 46.93870866666666
$ws.Cells.Item(16, 8).Value = 140.816126
$ws.Cells.Item(16, 9).Value = 0.117086244332836
$ws.Cells.Item(16, 10).Value = 0.117086244332836
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 63.67711
$ws.Cells.Item(16, 14).Value = 191.03133
$ws.Cells.Item(16, 15).Value = 0.4250937452800914
$ws.Cells.Item(16, 16).Value = 0.4250937452800915
$ws.Cells.Item(16, 17).Value = 2988.921315025286
$ws.Cells.Item(16, 18).Value = 26900.29183522758
$ws.Cells.Item(16, 19).Value = 0.04977263012422511
$ws.Cells.Item(16, 20).Value = 0.04977263012422514
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 46.93870866666666
$ws.Cells.Item(17, 8).Value = 140.816126
$ws.Cells.Item(17, 9).Value = 0.117086244332836
$ws.Cells.Item(17, 10).Value = 0.117086244332836
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 17.96553866666667
$ws.Cells.Item(17, 14).Value = 53.896616
$ws.Cells.Item(17, 15).Value = 0.119933805378222
$ws.Cells.Item(17, 16).Value = 0.119933805378222
$ws.Cells.Item(17, 17).Value = 843.2791855144018
$ws.Cells.Item(17, 18).Value = 7589.512669629616
$ws.Cells.Item(17, 19).Value = 0.0140425988402813
